$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.759.01'
$ws.Range('E2').Value = '  +3.97%  '
$ws.Range('D3').Value = '3.636.12'
$ws.Range('E3').Value = '  +2.72%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '203.99'
$ws.Range('E5').Value = '  +9.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '574.82'
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('D7').Value = '3.629.79'
$ws.Range('E7').Value = '  +2.62%  '
$ws.Range('E8').Value = '  +2.85%  '
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.690'
$ws.Range('E10').Value = '  +4.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '61.86'
$ws.Range('E11').Value = '  +18.05%  '
$ws.Range('E12').Value = '  +6.24%  '
$ws.Range('E13').Value = '  +13.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.19'
$ws.Range('E14').Value = '  +5.23%  '
$ws.Range('D15').Value = '4.211.81'
$ws.Range('E15').Value = '  +2.53%  '
$ws.Range('D16').Value = '3.632.76'
$ws.Range('E16').Value = '  +2.46%  '
$ws.Range('E17').Value = '  +1.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.09'
$ws.Range('E18').Value = '  +5.31%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.54'
$ws.Range('E19').Value = '  +4.18%  '
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '68.562.70'
$ws.Range('E20').Value = '  +3.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.09'
$ws.Range('E21').Value = '  +3.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '408.83'
$ws.Range('E22').Value = '  +4.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.98'
$ws.Range('E23').Value = '  +18.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.23'
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.26'
$ws.Range('E25').Value = '  +1.58%  '
$ws.Range('E26').Value = '  +3.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.74'
$ws.Range('E27').Value = '  +4.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.98'
$ws.Range('E28').Value = '  +14.35%  '
$ws.Range('E29').Value = '  +2.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.47'
$ws.Range('E30').Value = '  +7.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.08'
$ws.Range('E31').Value = '  +14.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.86'
$ws.Range('E32').Value = '  +3.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '675.02'
$ws.Range('E33').Value = '  +7.62%  '
$ws.Range('E34').Value = '  +2.73%  '
$ws.Range('E35').Value = '  +3.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '63.91'
$ws.Range('E36').Value = '  +1.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.39'
$ws.Range('E37').Value = '  +3.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.426'
$ws.Range('E38').Value = '  +8.49%  '
$ws.Range('D39').Value = '0.0₃0812'
$ws.Range('E39').Value = '  +8.26%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.28'
$ws.Range('E41').Value = '  +18.30%  '
$ws.Range('E42').Value = '  +5.52%  '
$ws.Range('D43').Value = '3.210.29'
$ws.Range('E43').Value = '  +8.74%  '
$ws.Range('E44').Value = '  +11.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('E46').Value = '  +28.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.93'
$ws.Range('E47').Value = '  +17.82%  '
$ws.Range('E48').Value = '  +5.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.96'
$ws.Range('E49').Value = '  +7.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.133'
$ws.Range('E50').Value = '  +2.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.11'
$ws.Range('E51').Value = '  +0.45%  '
